# Q3 Update - 2025
# UN-SAL.xlsx worksheet update:
#  - Refresh the "short-url" value (col B) shared across all data rows
#  - Insert a new "Haiti" country row (in alphabetical order) for year 2024
#  - Refresh refugee/asylum-seeker figures for several countries in the 2024 block
#  - Refresh the El Salvador 2024 "ooc"/"hst" totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newShortUrl = "8DNmNZ"

# The "short-url" value in column B is identical for every data row (2-157).
# Refresh it in one shot before the row insert below shifts everything down.
$ws.Range("B2:B157").Value = $newShortUrl

# Insert a new row for Haiti, in its correct alphabetical position,
# pushing Honduras..Venezuela down by one row (151 -> 158).
$ws.Rows.Item(151).Insert()

# Populate the new Haiti row (151)
$ws.Range("A151").Value = "1"
$ws.Range("B151").Value = "8DNmNZ"
$ws.Range("C151").Value = "1"
$ws.Range("D151").Value = "150"
$ws.Range("E151").Value = "2024"
$ws.Range("F151").Value = "81"
$ws.Range("G151").Value = "Haiti"
$ws.Range("H151").Value = "HAI"
$ws.Range("I151").Value = "HTI"
$ws.Range("J151").Value = "162"
$ws.Range("K151").Value = "El Salvador"
$ws.Range("L151").Value = "SAL"
$ws.Range("M151").Value = "SLV"
$ws.Range("N151").Value = "0"
$ws.Range("O151").Value = "5"
$ws.Range("P151").Value = "0"
$ws.Range("Q151").Value = "0"
$ws.Range("R151").Value = "0"
$ws.Range("S151").Value = "0"
$ws.Range("T151").Value = "0"
$ws.Range("U151").Value = "-"
$ws.Range("V151").Value = "0"

# Apply the updated figures (and renumbered "items" counters) for the
# rows that shifted down / changed values in the 2024 data block.
# Row 146
$ws.Range("O146").Value = "14"
# Row 147
$ws.Range("O147").Value = "8"
# Row 148
$ws.Range("N148").Value = "6"
$ws.Range("O148").Value = "55"
# Row 149
$ws.Range("O149").Value = "28"
# Row 150
$ws.Range("N150").Value = "8"
$ws.Range("O150").Value = "10"
# Row 152
$ws.Range("D152").Value = "151"
$ws.Range("N152").Value = "34"
$ws.Range("O152").Value = "44"
# Row 153
$ws.Range("D153").Value = "152"
# Row 154
$ws.Range("D154").Value = "153"
# Row 155
$ws.Range("D155").Value = "154"
$ws.Range("N155").Value = "33"
$ws.Range("O155").Value = "67"
# Row 156
$ws.Range("D156").Value = "155"
$ws.Range("T156").Value = "114393"
$ws.Range("V156").Value = "255433"
# Row 157
$ws.Range("D157").Value = "156"
$ws.Range("O157").Value = "6"
# Row 158
$ws.Range("D158").Value = "157"
$ws.Range("O158").Value = "30"
